$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("B2").Value = 20
$ws.Range("C3").Value = 50
$ws.Range("E3").Value = 81
$ws.Range("G3").Value = 69
$ws.Range("I3").Value = 102
$ws.Range("C6").Value = 259
$ws.Range("E6").Value = 233
$ws.Range("H6").Value = 228
$ws.Range("I6").Value = 293
$ws.Range("B7").Value = 280
$ws.Range("C7").Value = 350
$ws.Range("E7").Value = 360
$ws.Range("G7").Value = 388
$ws.Range("H7").Value = 358
$ws.Range("I7").Value = 469

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 29

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("H5").Value = 2
$ws.Range("H6").Value = 2

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 7

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("C3").Value = 3
$ws.Range("G3").Value = 3
$ws.Range("C5").Value = 26
$ws.Range("I5").Value = 11
$ws.Range("C6").Value = 29
$ws.Range("G6").Value = 27
$ws.Range("I6").Value = 26

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("H5").Value = 2
$ws.Range("C8").Value = 27
$ws.Range("H26").Value = 9
$ws.Range("C27").Value = 29
$ws.Range("G27").Value = 27
$ws.Range("I27").Value = 26
$ws.Range("I31").Value = 29
$ws.Range("I41").Value = 2
$ws.Range("E46").Value = 7
$ws.Range("I49").Value = 7
$ws.Range("E52").Value = 50
$ws.Range("I60").Value = 2
$ws.Range("E75").Value = 11
$ws.Range("I75").Value = 9
$ws.Range("B96").Value = 3
$ws.Range("B97").Value = 280
$ws.Range("C97").Value = 350
$ws.Range("E97").Value = 360
$ws.Range("G97").Value = 388
$ws.Range("H97").Value = 358
$ws.Range("I97").Value = 469

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("E6").Value = 39
$ws.Range("E7").Value = 50

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("E3").Value = 7
$ws.Range("I5").Value = 8
$ws.Range("E6").Value = 11
$ws.Range("I6").Value = 9

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("G3").Value = 1
$ws.Range("G5").Value = 2

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("H4").Value = 9
$ws.Range("H5").Value = 9

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("B2").Value = 2
$ws.Range("B6").Value = 3

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 7

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("C5").Value = 19
$ws.Range("C6").Value = 27
